# "Created TableRow and TableCell components"
#
# This workbook documents (on the "Props" sheet) which internal component
# each prop is threaded through on its way down the tree, and (on the
# "Tree" sheet) draws the component hierarchy itself. Introducing the new
# TableRow/TableCell layer means a couple of props now pass through an
# extra "5 - BodyContainer" hop, the "rows" prop stops being drilled down
# explicitly (it's consumed where the new components are created) and the
# ref that used to be called "tbodyRef" is renamed "tableBodyRef" to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Props")

# Row 15 - onItemsOpen: now also explicitly routed through the
# ScrollingContainer before reaching the BodyContainer.
$ws.Range("C15").Value = "3 - ScrollingContainer"
$ws.Range("D15").Value = "5 - BodyContainer"

# Row 24 - bodyContainerRef: now created on the BodyContainer itself
# (highlighted, since this is the newly introduced component) before
# being threaded on to SelectionRect and TableBody.
$ws.Range("D24").Value = "5 - BodyContainer"
$ws.Range("D24").Style = "Neutral"
$ws.Range("E24").Value = "6 - SelectionRect"
$ws.Range("F24").Value = "6 - TableBody"

# Row 25 - tbodyRef -> tableBodyRef: destination changes from TableBody to
# the new BodyContainer (highlighted).
$ws.Range("A25").Value = "tableBodyRef"
$ws.Range("C25").Value = "5 - BodyContainer"
$ws.Range("C25").Style = "Neutral"

# Row 28 ("rows") is removed - the rows data is now consumed internally by
# the new TableRow/TableCell components rather than drilled down further.
# The old row 29 ("columnResizeStart") shifts up into its place, leaving
# row 29 blank.
$ws.Range("A28").Value = "columnResizeStart"
$ws.Range("B28").Value = "4 - ResizingContainer"
$ws.Range("C28").Value = "6 - TableHead"
$ws.Range("A29:C29").Clear()
